# Apply the edit described by the diff:
# - Cell A7 text corrected from "Primary School student" to "Primary School Student"
# - Active cell selection moved to G12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of "Student" in A7
$ws.Range("A7").Value = "Primary School Student"

# Update the current selection to match the saved view state
$ws.Range("G12").Select()

$wb.Save()
